$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1505.9
$ws.Range("J32").Value = 1382.375
$ws.Range("L32").Value = 1382.375
$ws.Range("N32").Value = -2034.375
$ws.Range("H33").Value = 2466
$ws.Range("I33").Value = 347.63635
$ws.Range("K33").Value = 347.63635
$ws.Range("M33").Value = -118.63635
$ws.Range("H43").Value = 2599.8572
$ws.Range("I43").Value = 2405.7058
$ws.Range("J43").Value = 3425
$ws.Range("K43").Value = 2405.7058
$ws.Range("L43").Value = 3425
$ws.Range("M43").Value = -2336.7058
$ws.Range("N43").Value = -3563
$ws.Range("H100").Value = 3269.8
$ws.Range("I100").Value = 2837.5
$ws.Range("J100").Value = 4999
$ws.Range("K100").Value = 2837.5
$ws.Range("L100").Value = 4999
$ws.Range("M100").Value = -2296.5
$ws.Range("N100").Value = -6081
$ws.Range("H115").Value = 471.75
$ws.Range("I115").Value = 471.75
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1415.25
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 151.75
$ws.Range("N115").ClearContents()
$ws.Range("H129").Value = 761.8125
$ws.Range("I129").Value = 612.6
$ws.Range("J129").Value = 3000
$ws.Range("K129").Value = 1837.8
$ws.Range("L129").Value = 9000
$ws.Range("M129").Value = 3162.2
$ws.Range("N129").Value = -19000
$ws.Range("H132").Value = 2001.3673
$ws.Range("I132").Value = 2001.4681
$ws.Range("K132").Value = 6004.4043
$ws.Range("M132").Value = -3474.4043
$ws.Range("H138").Value = 4621.413
$ws.Range("J138").Value = 4416.778
$ws.Range("L138").Value = 13250.334
$ws.Range("N138").Value = -23530.334
$ws.Range("H141").Value = 2184.5715
$ws.Range("I141").Value = 2184.5715
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6553.7145
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1373.7145
$ws.Range("N141").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8516.382
$ws.Range("I32").Value = 4718.871
$ws.Range("K32").Value = 4718.871
$ws.Range("M32").Value = -4431.871
$ws.Range("H74").Value = 5207.1577
$ws.Range("I74").Value = 3077.8
$ws.Range("K74").Value = 3077.8
$ws.Range("M74").Value = -2203.8
$ws.Range("H77").Value = 5207.1577
$ws.Range("I77").Value = 3077.8
$ws.Range("K77").Value = 15389
$ws.Range("M77").Value = -11021
$ws.Range("H110").Value = 1605.909
$ws.Range("I110").Value = 1412
$ws.Range("J110").Value = 2265.2
$ws.Range("K110").Value = 1412
$ws.Range("L110").Value = 2265.2
$ws.Range("M110").Value = 633
$ws.Range("N110").Value = -6355.2

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 35000
$ws.Range("L21").Value = 35000
$ws.Range("N21").Value = -35472
$ws.Range("H22").Value = 701
$ws.Range("I22").Value = 510.35294
$ws.Range("K22").Value = 510.35294
$ws.Range("M22").Value = -337.35294
$ws.Range("H105").Value = 2227544
$ws.Range("I105").Value = 2671185.5
$ws.Range("J105").Value = 9335.666999999999
$ws.Range("K105").Value = 2671185.5
$ws.Range("L105").Value = 9335.666999999999
$ws.Range("M105").Value = -2669438.5
$ws.Range("N105").Value = -12829.667

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 577
$ws.Range("I16").Value = 493
$ws.Range("J16").Value = 703
$ws.Range("K16").Value = 493
$ws.Range("L16").Value = 703
$ws.Range("M16").Value = -206
$ws.Range("N16").Value = -1277
$ws.Range("H86").Value = 5719650.5
$ws.Range("I86").Value = 13337752
$ws.Range("J86").Value = 6074.25
$ws.Range("K86").Value = 13337752
$ws.Range("L86").Value = 6074.25
$ws.Range("M86").Value = -13336629
$ws.Range("N86").Value = -8320.25
$ws.Range("H89").Value = 5719650.5
$ws.Range("I89").Value = 13337752
$ws.Range("J89").Value = 6074.25
$ws.Range("K89").Value = 66688760
$ws.Range("L89").Value = 30371.25
$ws.Range("M89").Value = -66683144
$ws.Range("N89").Value = -41603.25
$ws.Range("H96").Value = 80949.75
$ws.Range("J96").Value = 80949.75
$ws.Range("L96").Value = 80949.75
$ws.Range("N96").Value = -86441.75
$ws.Range("H99").Value = 4588.4043
$ws.Range("I99").Value = 5780.6924
$ws.Range("J99").Value = 3112.238
$ws.Range("K99").Value = 5780.6924
$ws.Range("L99").Value = 3112.238
$ws.Range("M99").Value = -4282.6924
$ws.Range("N99").Value = -6108.237999999999
$ws.Range("H113").Value = 577
$ws.Range("I113").Value = 493
$ws.Range("J113").Value = 703
$ws.Range("K113").Value = 493
$ws.Range("L113").Value = 703
$ws.Range("M113").Value = 1677
$ws.Range("N113").Value = -5043
$ws.Range("H126").Value = 4588.4043
$ws.Range("I126").Value = 5780.6924
$ws.Range("J126").Value = 3112.238
$ws.Range("K126").Value = 17342.0772
$ws.Range("L126").Value = 9336.714
$ws.Range("M126").Value = -14872.0772
$ws.Range("N126").Value = -14276.714
$ws.Range("H132").Value = 4919.5435
$ws.Range("I132").Value = 3151.3076
$ws.Range("J132").Value = 14771.143
$ws.Range("K132").Value = 9453.9228
$ws.Range("L132").Value = 44313.429
$ws.Range("M132").Value = -6923.9228
$ws.Range("N132").Value = -49373.429

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 57.166668
$ws.Range("J2").Value = 36.57143
$ws.Range("L2").Value = 219.42858
$ws.Range("N2").Value = -445.42858
$ws.Range("H5").Value = 1343.0714
$ws.Range("I5").Value = 873.8570999999999
$ws.Range("K5").Value = 2621.5713
$ws.Range("M5").Value = -2509.5713
$ws.Range("H38").Value = 13.333333
$ws.Range("I38").Value = 11.6
$ws.Range("K38").Value = 34.8
$ws.Range("M38").Value = 312.2
$ws.Range("H86").Value = 489.5
$ws.Range("J86").Value = 489.5
$ws.Range("L86").Value = 1468.5
$ws.Range("N86").Value = -3840.5
$ws.Range("H89").Value = 489.5
$ws.Range("J89").Value = 489.5
$ws.Range("L89").Value = 4405.5
$ws.Range("N89").Value = -16261.5
$ws.Range("H92").Value = 462.4762
$ws.Range("J92").Value = 469.41177
$ws.Range("L92").Value = 1408.23531
$ws.Range("N92").Value = -3904.23531
$ws.Range("H112").Value = 18000
$ws.Range("I112").Value = 10000
$ws.Range("K112").Value = 30000
$ws.Range("M112").Value = -28892
$ws.Range("H114").Value = 429.25
$ws.Range("I114").Value = 393.2
$ws.Range("K114").Value = 1179.6
$ws.Range("M114").Value = 2074.4
$ws.Range("H135").Value = 1343.0714
$ws.Range("I135").Value = 873.8570999999999
$ws.Range("K135").Value = 7864.7139
$ws.Range("M135").Value = -5329.7139

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 4889.625
$ws.Range("I113").Value = 4889.625
$ws.Range("K113").Value = 4889.625
$ws.Range("M113").Value = -2719.625

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8057.206
$ws.Range("J2").Value = 14777.723
$ws.Range("L2").Value = 14777.723
$ws.Range("N2").Value = -15001.723
$ws.Range("H7").Value = 4124.636
$ws.Range("I7").Value = 4207.1
$ws.Range("K7").Value = 4207.1
$ws.Range("M7").Value = -4095.1
$ws.Range("H40").Value = 9873.333000000001
$ws.Range("I40").Value = 7676.75
$ws.Range("J40").Value = 10971.625
$ws.Range("K40").Value = 7676.75
$ws.Range("L40").Value = 10971.625
$ws.Range("M40").Value = -7540.75
$ws.Range("N40").Value = -11243.625
$ws.Range("H55").Value = 479.41666
$ws.Range("I55").Value = 315.64285
$ws.Range("K55").Value = 315.64285
$ws.Range("M55").Value = -142.64285
$ws.Range("H61").Value = 4860.7417
$ws.Range("I61").Value = 1909.7368
$ws.Range("K61").Value = 1909.7368
$ws.Range("M61").Value = -1707.7368
$ws.Range("H113").Value = 4860.7417
$ws.Range("I113").Value = 1909.7368
$ws.Range("K113").Value = 1909.7368
$ws.Range("M113").Value = 260.2632000000001
$ws.Range("H126").Value = 4124.636
$ws.Range("I126").Value = 4207.1
$ws.Range("K126").Value = 12621.3
$ws.Range("M126").Value = -10151.3

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3644.88
$ws.Range("I96").Value = 1525.0555
$ws.Range("J96").Value = 9095.857
$ws.Range("K96").Value = 1525.0555
$ws.Range("L96").Value = 9095.857
$ws.Range("M96").Value = -152.0554999999999
$ws.Range("N96").Value = -11841.857
$ws.Range("H113").Value = 249.65384
$ws.Range("I113").Value = 185.57143
$ws.Range("K113").Value = 556.71429
$ws.Range("M113").Value = 1613.28571
$ws.Range("H126").Value = 1411
$ws.Range("I126").Value = 1411
$ws.Range("K126").Value = 4233
$ws.Range("M126").Value = -1763
